$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-7 (columns A-E)
$data = @(
    @(0, "Skanda-4865A-654.50-obverse.jpg", "Gupta/Skandagupta", "archer type", "King standing and holding an Indian long bow in his left hand and an arrow in his right Garuda standard at left;Brāhmī legend under arm: skanda"),
    @(1, "Skanda-4865A-654.50-reverse.jpg", "Gupta/Skandagupta", "archer type", "Lakshmi seated facing, holding long-stemmed lotus and diadem;Brāhmī legend at right: sri skandaguptah;circular Brāhmī legend around"),
    @(2, "Skanda-4865-148.06-obverse.jpg", "Gupta/Skandagupta", "archer type", "King standing and holding an Indian long bow in his left hand and an arrow in his right Garuda standard at left;Brāhmī legend under arm: skanda;circular Brāhmī legend around"),
    @(3, "Skanda-4865-148.06-reverse.jpg", "Gupta/Skandagupta", "archer type", "Lakshmi seated facing, holding long-stemmed lotus and diadem;Brāhmī legend at right: kramadityah"),
    @(4, "Skanda-4866-421.04-obverse.jpg", "Gupta/Skandagupta", "archer type;King and Lakshmi type", "King standing at left, facing right and holding an Indian long bow in his left hand and an arrow in his right hand, Lakshmi standing at right, facing left and offering an indistinct object to the king Garuda standard between the two figures,"),
    @(5, "Skanda-4866-421.04-reverse.jpg", "Gupta/Skandagupta", "archer type", "Lakshmi seated facing, holding long-stemmed lotus and diadem;Brāhmī legend at right: sri skandaguptah;circular Brāhmī legend around")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
